$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.582.72"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.15%  "

$ws.Range("D3").Value = "'3.909.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +3.08%  "

$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "'603.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").Value = "'164.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.79%  "

$ws.Range("D7").Value = "'3.909.51"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.15%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("D9").Value = "'0.527"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.31%  "

$ws.Range("D10").Value = "'0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.27%  "

$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("E12").Value = "  -0.68%  "

$ws.Range("E13").Value = "  -1.04%  "

$ws.Range("D14").Value = "'36.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.53%  "

$ws.Range("D15").Value = "'4.558.21"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.00%  "

$ws.Range("D16").Value = "'3.911.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.19%  "

$ws.Range("D17").Value = "'68.705.82"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.10%  "

$ws.Range("D18").Value = "'7.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.61%  "

$ws.Range("E19").Value = "  -1.62%  "

$ws.Range("D20").Value = "'17.02"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.74%  "

$ws.Range("D21").Value = "'11.18"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.00%  "

$ws.Range("D22").Value = "'484.25"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.97%  "

$ws.Range("D23").Value = "'0.718"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.40%  "

$ws.Range("E24").Value = "  +10.96%  "

$ws.Range("D25").Value = "'84.36"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.70%  "

$ws.Range("E26").Value = "  -2.20%  "

$ws.Range("D27").Value = "'11.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.76%  "

$ws.Range("D28").Value = "'10.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.87%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  -2.02%  "

$ws.Range("D31").Value = "'4.056.04"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.25%  "

$ws.Range("D32").Value = "'7.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.16%  "

$ws.Range("D33").Value = "'2.37"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.03%  "

$ws.Range("D34").Value = "'31.98"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.18%  "

$ws.Range("D35").Value = "'3.848.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.97%  "

$ws.Range("E36").Value = "  -1.94%  "

$ws.Range("D37").Value = "'1.04"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.55%  "

$ws.Range("E38").Value = "  -0.28%  "

$ws.Range("D39").Value = "'5.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.82%  "

$ws.Range("E40").Value = "  +0.04%  "

$ws.Range("D41").Value = "'3.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("E42").Value = "  -3.27%  "

$ws.Range("D43").Value = "'431.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.19%  "

$ws.Range("E44").Value = "  -0.10%  "

$ws.Range("E45").Value = "  -1.12%  "

$ws.Range("D47").Value = "'8.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.69%  "

$ws.Range("D48").Value = "'26.53"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.76%  "

$ws.Range("D49").Value = "'141.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").Value = "'2.814.97"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("E51").Value = "  -2.57%  "

